# TodoList.xlsx edit:
#  - "exc add chapter annoucments" (G6) trimmed to "add chapter annoucments"
#  - "exec delete chapter annoouncements" (G7) trimmed to " delete chapter annoouncements"
#  - new cell G8 gets "exec only" (the "exec" piece split out of the two items above)
#  - A5, D5, G6, G7 get a red highlight fill to flag these TODO rows
#  - selection moved to C13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- text edits: split "exec ..." prefixes out of the two announcement items ---
$ws.Range("G6").Value = "add chapter annoucments"
$ws.Range("G7").Value = " delete chapter annoouncements"
$ws.Range("G8").Value = "exec only"

# --- highlight formatting (red fill) ---
$ws.Range("A5").Interior.Color = 255
$ws.Range("D5").Interior.Color = 255
$ws.Range("G6").Interior.Color = 255
$ws.Range("G7").Interior.Color = 255

# --- move the active selection to C13 ---
$ws.Range("C13").Select()
